$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 19, shifting rows 19:25 down to 20:26
$ws.Rows.Item(19).Insert(-4121)

# Fill in the new row 19 with the inserted record's data
$ws.Cells.Item(19, 1).Value = 11
$ws.Cells.Item(19, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(19, 3).Value = "Bíobío"
$ws.Cells.Item(19, 4).Value = 45079
$ws.Cells.Item(19, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(19, 5).Value = 8
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100104
$ws.Cells.Item(19, 8).Value = "Frutos de pepita"
$ws.Cells.Item(19, 9).Value = 100104003
$ws.Cells.Item(19, 10).Value = "Membrillo"
$ws.Cells.Item(19, 11).Value = "Champion"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 270
$ws.Cells.Item(19, 14).Value = 11000
$ws.Cells.Item(19, 15).Value = 12000
$ws.Cells.Item(19, 16).Value = 11444
$ws.Cells.Item(19, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(19, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(19, 19).Value = 636
$ws.Cells.Item(19, 20).Value = 18
